$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47 - pushes existing rows 47:109 down to 48:110
$ws.Rows("47").Insert()

# Populate the newly inserted row 47 with the latest weekly price record
$ws.Range("A47").Value = 10
$ws.Range("B47").Value = "Vega Modelo de Temuco"
$ws.Range("C47").Value = "La Araucanía"
$ws.Range("D47").Value = 44533
$ws.Range("E47").Value = 9
$ws.Range("F47").Value = 100112012
$ws.Range("G47").Value = "Espinaca"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 110
$ws.Range("K47").Value = 8000
$ws.Range("L47").Value = 8000
$ws.Range("M47").Value = 8000
$ws.Range("N47").Value = "$/docena de atados"
$ws.Range("O47").Value = "Región de La Araucanía"
$ws.Range("P47").Value = 2667
$ws.Range("Q47").Value = 3
$ws.Range("R47").Value = "Hortaliza"
